{"js": "// Replace each three-digit-by-one-digit multiplication expression in the\n// practice-sheet table with its new value, matching the author's commit.\n// Every \"old\" string below occurs exactly once in the document, so a plain\n// literal (non-wildcard) search-and-replace is safe and unambiguous.\nconst replacements = [\n  [\"962\u00d79=\", \"534\u00d74=\"],\n  [\"609\u00d73=\", \"528\u00d75=\"],\n  [\"544\u00d78=\", \"844\u00d72=\"],\n  [\"605\u00d79=\", \"751\u00d77=\"],\n  [\"994\u00d73=\", \"316\u00d75=\"],\n  [\"983\u00d76=\", \"678\u00d79=\"],\n  [\"400\u00d76=\", \"845\u00d77=\"],\n  [\"859\u00d75=\", \"823\u00d77=\"],\n  [\"275\u00d79=\", \"950\u00d76=\"],\n  [\"929\u00d79=\", \"281\u00d76=\"],\n  [\"594\u00d73=\", \"298\u00d78=\"],\n  [\"217\u00d74=\", \"642\u00d76=\"],\n  [\"863\u00d76=\", \"757\u00d78=\"],\n  [\"340\u00d72=\", \"241\u00d76=\"],\n  [\"145\u00d78=\", \"182\u00d72=\"],\n  [\"453\u00d77=\", \"842\u00d77=\"],\n  [\"919\u00d75=\", \"896\u00d74=\"],\n  [\"132\u00d73=\", \"944\u00d73=\"],\n  [\"215\u00d79=\", \"577\u00d76=\"],\n  [\"294\u00d73=\", \"227\u00d76=\"],\n  [\"866\u00d73=\", \"757\u00d76=\"],\n  [\"542\u00d73=\", \"920\u00d75=\"],\n  [\"658\u00d78=\", \"776\u00d75=\"],\n  [\"321\u00d76=\", \"267\u00d72=\"],\n  [\"667\u00d79=\", \"572\u00d74=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const r of results.items) {\n    r.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each three-digit-by-one-digit multiplication expression in the\n# practice-sheet table with its new value, matching the author's commit.\n# Every \"old\" string occurs exactly once in the document, so a literal\n# (non-wildcard) Find/Replace is safe and unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"962\u00d79=\", \"534\u00d74=\"),\n    @(\"609\u00d73=\", \"528\u00d75=\"),\n    @(\"544\u00d78=\", \"844\u00d72=\"),\n    @(\"605\u00d79=\", \"751\u00d77=\"),\n    @(\"994\u00d73=\", \"316\u00d75=\"),\n    @(\"983\u00d76=\", \"678\u00d79=\"),\n    @(\"400\u00d76=\", \"845\u00d77=\"),\n    @(\"859\u00d75=\", \"823\u00d77=\"),\n    @(\"275\u00d79=\", \"950\u00d76=\"),\n    @(\"929\u00d79=\", \"281\u00d76=\"),\n    @(\"594\u00d73=\", \"298\u00d78=\"),\n    @(\"217\u00d74=\", \"642\u00d76=\"),\n    @(\"863\u00d76=\", \"757\u00d78=\"),\n    @(\"340\u00d72=\", \"241\u00d76=\"),\n    @(\"145\u00d78=\", \"182\u00d72=\"),\n    @(\"453\u00d77=\", \"842\u00d77=\"),\n    @(\"919\u00d75=\", \"896\u00d74=\"),\n    @(\"132\u00d73=\", \"944\u00d73=\"),\n    @(\"215\u00d79=\", \"577\u00d76=\"),\n    @(\"294\u00d73=\", \"227\u00d76=\"),\n    @(\"866\u00d73=\", \"757\u00d76=\"),\n    @(\"542\u00d73=\", \"920\u00d75=\"),\n    @(\"658\u00d78=\", \"776\u00d75=\"),\n    @(\"321\u00d76=\", \"267\u00d72=\"),\n    @(\"667\u00d79=\", \"572\u00d74=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
